# The commit removes the stray word " den" from the sentence
#   "...i hovedutvalg for kultur, idrett og folkehelse den {moetedato}."
# so that it reads
#   "...i hovedutvalg for kultur, idrett og folkehelse {moetedato}."
#
# (the rest of the upstream diff is Word's background spell/grammar
# checker re-flowing the surrounding runs and stamping <w:proofErr/>
# markers + bumping schema namespaces on save - it does not change any
# visible text, so there is nothing else to "edit" here.)

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute(" den")
if ($found) {
    $rng.Delete()
}
